$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("title21", "No Stamina Cost", "スタミナ消費なし", "不消耗体力"),
    @("toggle55", "Enable No Stamina Cost", "スタミナ消費なしを有効化", "启用不消耗体力"),
    @("tooltip21", "Enable or disable no stamina cost while fishing.", "釣り中のスタミナ消費なしを有効または無効にします。", "启用或禁用钓鱼时不消耗体力。")
)

$row = 108
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}
